# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the last row
# (bf81812a-8ebd-4342-991e-cfeffb0c96aa) in both the zh-cn and de-de sheets,
# as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-20 04:35:10"
$wsZhCn.Range("H4").Value = "2016-03-20 04:35:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-20 04:35:13"
$wsDeDe.Range("H4").Value = "2016-03-20 04:35:53"
